$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.420.73'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.636.21'
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.17'
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.532'
$ws.Range("E6").Value = '  +4.28%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.98'
$ws.Range("E8").Value = '  -4.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.255'
$ws.Range("E9").Value = '  -2.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0608'
$ws.Range("E10").Value = '  -1.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0885'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.868.62'
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.634.67'
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.578'
$ws.Range("E14").Value = '  +2.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.00'
$ws.Range("E15").Value = '  -2.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.09'
$ws.Range("E16").Value = '  -2.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.417.67'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.70'
$ws.Range("E18").Value = '  -3.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0721'
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.50'
$ws.Range("E20").Value = '  -0.88%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").Value = '  -3.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.63'
$ws.Range("E23").Value = '  +3.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.97'
$ws.Range("E24").Value = '  -3.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.09'
$ws.Range("E25").Value = '  +2.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.96'
$ws.Range("E26").Value = '  -3.02%  '
$ws.Range("E27").Value = '  +1.33%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.50'
$ws.Range("E29").Value = '  -3.68%  '
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0486'
$ws.Range("E31").Value = '  -2.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.28'
$ws.Range("E32").Value = '  -0.78%  '
$ws.Range("E33").Value = '  +3.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.408.50'
$ws.Range("E34").Value = '  -2.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").Value = '  +1.73%  '
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("E37").Value = '  -0.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.874'
$ws.Range("E38").Value = '  -4.30%  '
$ws.Range("E39").Value = '  -2.05%  '
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.814'
$ws.Range("E42").Value = '  +3.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.47'
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.61'
$ws.Range("E45").Value = '  -2.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.777.00'
$ws.Range("E47").Value = '  -3.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.63'
$ws.Range("E48").Value = '  -3.29%  '
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0990'
$ws.Range("E50").Value = '  -1.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.68'
$ws.Range("E51").Value = '  -1.46%  '
